$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04408699999999999
$ws.Range("H2").Value = 0.132261
$ws.Range("I2").Value = 0.007006504090795892
$ws.Range("J2").Value = 0.007006504090795892
$ws.Range("Q2").Value = 0.04810291422133332
$ws.Range("R2").Value = 0.432926227992
$ws.Range("S2").Value = 0.006690321216221995
$ws.Range("T2").Value = 0.006690321216221995

$ws.Range("G3").Value = 0.04408699999999999
$ws.Range("H3").Value = 0.132261
$ws.Range("I3").Value = 0.007006504090795892
$ws.Range("J3").Value = 0.007006504090795892
$ws.Range("Q3").Value = 0.002273331459333333
$ws.Range("R3").Value = 0.020459983134
$ws.Range("S3").Value = 0.0003161828745738959
$ws.Range("T3").Value = 0.0003161828745738959

$ws.Range("I4").Value = 0.9567202519440571
$ws.Range("J4").Value = 0.9567202519440571
$ws.Range("S4").Value = 0.9135462873673272
$ws.Range("T4").Value = 0.9135462873673272

$ws.Range("I5").Value = 0.9567202519440571
$ws.Range("J5").Value = 0.9567202519440571
$ws.Range("S5").Value = 0.04317396457672975
$ws.Range("T5").Value = 0.04317396457672974

$ws.Range("I6").Value = 0.03627324396514701
$ws.Range("J6").Value = 0.03627324396514701
$ws.Range("S6").Value = 0.03463633939784837
$ws.Range("T6").Value = 0.03463633939784837

$ws.Range("I7").Value = 0.03627324396514701
$ws.Range("J7").Value = 0.03627324396514701
$ws.Range("S7").Value = 0.001636904567298641
$ws.Range("T7").Value = 0.00163690456729864
